$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# Shapes "Google Shape;189;p5" .. "Google Shape;200;p5" are the small
# rounded-rectangle labels in the process diagram. Each currently has a
# text inset (108000 or 72000 EMU converted to points) and shape-to-fit-text
# autosizing (a:spAutoFit). Switch them to zero insets and no autofit
# (a:noAutofit), matching the target layout.
$targetNames = @(
    "Google Shape;189;p5",
    "Google Shape;190;p5",
    "Google Shape;191;p5",
    "Google Shape;192;p5",
    "Google Shape;193;p5",
    "Google Shape;194;p5",
    "Google Shape;195;p5",
    "Google Shape;196;p5",
    "Google Shape;197;p5",
    "Google Shape;198;p5",
    "Google Shape;199;p5",
    "Google Shape;200;p5"
)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($targetNames -contains $shp.Name) {
        $tf = $shp.TextFrame
        $tf.MarginLeft = 0
        $tf.MarginRight = 0
        $tf.MarginTop = 0
        $tf.MarginBottom = 0
        $tf.AutoSize = 0
    }
}
